$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the existing (invalid / test) data rows, keep only valid rows
$ws.Range("A1:C7").ClearContents()

# Header row
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Id"

# Valid data row
$ws.Range("A2").Value = "Guyal15"
$ws.Range("B2").Value = "123456a!"
$ws.Range("C2").Value = 316471465

# Column widths to roughly match target layout
$ws.Columns.Item(1).ColumnWidth = 9.73
$ws.Columns.Item(2).ColumnWidth = 9.73
$ws.Columns.Item(3).ColumnWidth = 9.31

# Selection matches the saved state in the target file
$ws.Range("C1").Select() | Out-Null
